# This script rotates the data of rows 6, 7 and 8 (species observation
# records) one step "up": the old row 7 values move into row 6, the old
# row 8 values move into row 7, and the old row 6 values move into row 8.
# The K/L/M/N/AC "activity / comment" block that used to sit on row 6
# travels along with that record, ending up on row 8; row 6 and 7 end up
# without that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that participate in the rotation (values differ row to row).
$cols = @("A","B","D","E","F","G","H","Q","R")

# Capture current values of rows 6, 7, 8 for those columns.
$row6 = @{}
$row7 = @{}
$row8 = @{}
foreach ($col in $cols) {
    $row6[$col] = $ws.Range("$col`6").Value2
    $row7[$col] = $ws.Range("$col`7").Value2
    $row8[$col] = $ws.Range("$col`8").Value2
}

# Capture the extra "activity" block (K, L, M, N) and the public comment
# (AC) that currently live on row 6; they move together with that record.
$k6 = $ws.Range("K6").Value2
$l6 = $ws.Range("L6").Value2
$m6 = $ws.Range("M6").Value2
$n6 = $ws.Range("N6").Value2
$ac6 = $ws.Range("AC6").Value2

# Write the rotated values: row6 <- row7, row7 <- row8, row8 <- row6.
foreach ($col in $cols) {
    $ws.Range("$col`6").Value = $row7[$col]
    $ws.Range("$col`7").Value = $row8[$col]
    $ws.Range("$col`8").Value = $row6[$col]
}

# Move the activity/comment block from row 6 to row 8, clearing it from
# rows 6 and 7.
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("AC6").ClearContents()

$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("AC7").ClearContents()

$ws.Range("K8").Value = $k6
$ws.Range("L8").Value = $l6
$ws.Range("M8").Value = $m6
$ws.Range("N8").Value = $n6
$ws.Range("AC8").Value = $ac6
